$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated forecast-error values after fixing the naive component forecaster bug.
$data = @(
    @(2, 0.1604947413933812,  0.671979494075615,  0.9800629839820976, 0.9899813048649442, 0.9866055495518523, 51),
    @(3, -0.01688632295344274, 0.678798267575447,  1.154921742939148,  1.074672853913761,  1.085449495453315,  50),
    @(4, 0.09175790792555,    0.6751575700206842, 1.246489116971403,  1.11646277007852,   1.124216491624296,  49),
    @(5, -0.03150304691393851, 0.6526513673133957, 1.099963246005613,  1.04879132624446,   1.059411705416517,  48),
    @(6, 0.09092125959470221, 0.6433332090867682, 1.118209597697466,  1.057454300524361,  1.06492822115538,   47),
    @(7, 0.06426616516546159, 0.6673738110489854, 1.263617861539698,  1.124107584504125,  1.137333693445019,  38),
    @(8, 0.07636997984574115, 0.6282867277610423, 1.249410503617163,  1.117770326863781,  1.130540563949073,  37),
    @(9, -0.1896082154670432, 0.6391424148142104, 0.8611962251771343, 0.9280065868177523, 0.9320294493904645, 20),
    @(10, 0.02217556234913506, 0.5877922150015374, 0.7922902390694693, 0.8901068694653858, 0.926165044936718,  13),
    @(11, -0.0636988451256798, 0.3782091602451033, 0.2022309754915021, 0.4497009845347262, 0.4977115537903916, 5)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
}
